# This workbook tracks weekly price observations for Cilantro at
# "Femacal de La Calera". A new weekly observation row is inserted at
# row 550, pushing all subsequent rows (550-639) down by one (to 551-640).
# We then populate the newly inserted row 550 with its data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 550; this shifts rows 550:639 down to 551:640
# (Excel copies the formatting of the row above into the new row, which is
# why the D column keeps its date style automatically).
$ws.Rows.Item(550).Insert()

# Populate the newly inserted row 550 with the new observation's data.
$ws.Range("A550").Value = 3
$ws.Range("B550").Value = "Femacal de La Calera"
$ws.Range("C550").Value = "Coquimbo"
$ws.Range("D550").Value = 45180
$ws.Range("E550").Value = 5
$ws.Range("F550").Value = 100112040
$ws.Range("G550").Value = "Cilantro"
$ws.Range("H550").Value = "Sin especificar"
$ws.Range("I550").Value = "Primera"
$ws.Range("J550").Value = 200
$ws.Range("K550").Value = 4000
$ws.Range("L550").Value = 4500
$ws.Range("M550").Value = 4300
$ws.Range("N550").Value = "$/docena de atados (3 kilos)"
$ws.Range("O550").Value = "Provincia de Quillota"
$ws.Range("P550").Value = 1433
$ws.Range("Q550").Value = 3
$ws.Range("R550").Value = "Hortaliza"
